# Update the "想去人数" (want-to-go count) figures for the two events
# that appear on both the "展览" sheet and the "全部类型" sheet.
#   F2: 367 -> 368
#   F3: 380 -> 381

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 368
    $ws.Range("F3").Value = 381
}
